$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '243.76'
Set-TextValue 'D3' '23.94'
$ws.Range('B4').Value = 'LEO'
$ws.Range('C4').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D4' '3.502'
$ws.Range('E4').Value = '3LEOLEO'
$ws.Range('B5').Value = 'HuobiToken'
$ws.Range('C5').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D5' '5.142'
$ws.Range('E5').Value = '4HuobiTokenHT'
$ws.Range('B6').Value = 'Cronos'
$ws.Range('C6').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D6' '0.05746'
$ws.Range('E6').Value = '5CronosCRO'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D7' '6.474'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D8' '3.145'
$ws.Range('E8').Value = '7GateTokenGT'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.8096'
$ws.Range('E9').Value = '8MXTokenMX'
$ws.Range('B10').Value = 'FTXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D10' '0.8345'
$ws.Range('E10').Value = '9FTXTokenFTT'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1337'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.06945'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D13' '0.03122'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D14' '0.02850'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D15' '0.09366'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D16' '3.747'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D17' '0.001513'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D18' '0.04676'
$ws.Range('E18').Value = '17CoinExTokenCET'
Set-TextValue 'D19' '0.006240'
Set-TextValue 'D20' '0.001239'
Set-TextValue 'D21' '0.004271'
Set-TextValue 'D22' '0.00008708'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D23' '2.082'
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D24' '0.009647'
$ws.Range('E24').Value = '23OneONEBestin24h'
Set-TextValue 'D26' '0.1339'
Set-TextValue 'D40' '0.03615'
Set-TextValue 'D41' '0.006446'
Set-TextValue 'D42' '0.1048'
$ws.Range('E43').Value = '42CEJICEJI'
Set-TextValue 'D44' '0.007372'
Set-TextValue 'D45' '0.00005309'
Set-TextValue 'D48' '0.002275'
Set-TextValue 'D50' '0.0002002'

Write-Output "Applied 83 cell updates"
